$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

Set-TextValue $ws.Range("D2") "26.768.40"
Set-TextValue $ws.Range("E2") "  +0.43%  "
Set-TextValue $ws.Range("D3") "1.651.36"
Set-TextValue $ws.Range("E3") "  +0.90%  "
Set-TextValue $ws.Range("E4") "  +0.05%  "
Set-TextValue $ws.Range("D5") "215.76"
Set-TextValue $ws.Range("E5") "  +1.35%  "
Set-TextValue $ws.Range("D6") "0.505"
Set-TextValue $ws.Range("E6") "  +0.80%  "
Set-TextValue $ws.Range("E7") "  -0.02%  "
Set-TextValue $ws.Range("D8") "0.254"
Set-TextValue $ws.Range("E8") "  -0.51%  "
Set-TextValue $ws.Range("D9") "0.0630"
Set-TextValue $ws.Range("E9") "  +0.95%  "
Set-TextValue $ws.Range("D10") "19.44"
Set-TextValue $ws.Range("E10") "  +1.80%  "
Set-TextValue $ws.Range("D11") "0.0847"
Set-TextValue $ws.Range("E11") "  +0.69%  "
Set-TextValue $ws.Range("D12") "1.879.29"
Set-TextValue $ws.Range("E12") "  +0.74%  "
Set-TextValue $ws.Range("D15") "0.536"
Set-TextValue $ws.Range("E15") "  +1.50%  "
Set-TextValue $ws.Range("D16") "66.63"
Set-TextValue $ws.Range("E16") "  +5.11%  "
Set-TextValue $ws.Range("D17") "26.795.84"
Set-TextValue $ws.Range("D18") "0.0₃0757"
Set-TextValue $ws.Range("E18") "  +1.68%  "
Set-TextValue $ws.Range("D19") "221.54"
Set-TextValue $ws.Range("E19") "  +1.61%  "
Set-TextValue $ws.Range("E20") "  +0.03%  "
Set-TextValue $ws.Range("D21") "4.42"
Set-TextValue $ws.Range("E21") "  +2.43%  "
Set-TextValue $ws.Range("D22") "6.41"
Set-TextValue $ws.Range("E22") "  +2.63%  "
Set-TextValue $ws.Range("D23") "9.59"
Set-TextValue $ws.Range("E23") "  +0.68%  "
Set-TextValue $ws.Range("D24") "2.17"
Set-TextValue $ws.Range("E24") "  +13.11%  "
Set-TextValue $ws.Range("D25") "147.97"
Set-TextValue $ws.Range("E25") "  -0.81%  "
Set-TextValue $ws.Range("E26") "  +0.04%  "
Set-TextValue $ws.Range("E27") "  -0.31%  "
Set-TextValue $ws.Range("D28") "7.12"
Set-TextValue $ws.Range("E28") "  +3.25%  "
Set-TextValue $ws.Range("D29") "15.99"
Set-TextValue $ws.Range("E29") "  +3.34%  "
Set-TextValue $ws.Range("D30") "0.0524"
Set-TextValue $ws.Range("E30") "  +0.65%  "
Set-TextValue $ws.Range("E31") "  +0.54%  "
Set-TextValue $ws.Range("D32") "3.45"
Set-TextValue $ws.Range("E32") "  +4.90%  "
Set-TextValue $ws.Range("D33") "3.07"
Set-TextValue $ws.Range("E33") "  +4.32%  "
Set-TextValue $ws.Range("E34") "  +4.89%  "
Set-TextValue $ws.Range("D35") "1.299.24"
Set-TextValue $ws.Range("E35") "  +9.78%  "
Set-TextValue $ws.Range("E36") "  +5.36%  "
Set-TextValue $ws.Range("E37") "  +0.93%  "
Set-TextValue $ws.Range("D38") "0.830"
Set-TextValue $ws.Range("E38") "  +2.52%  "
Set-TextValue $ws.Range("D39") "0.527"
Set-TextValue $ws.Range("E39") "  +3.71%  "
Set-TextValue $ws.Range("E40") "  +0.04%  "
Set-TextValue $ws.Range("D41") "0.818"
Set-TextValue $ws.Range("E41") "  +3.07%  "
Set-TextValue $ws.Range("E42") "  -2.41%  "
Set-TextValue $ws.Range("D43") "5.44"
Set-TextValue $ws.Range("E43") "  +0.36%  "
Set-TextValue $ws.Range("D44") "1.791.61"
Set-TextValue $ws.Range("E44") "  +1.07%  "
Set-TextValue $ws.Range("D45") "93.94"
Set-TextValue $ws.Range("E45") "  +1.39%  "
Set-TextValue $ws.Range("D46") "61.16"
Set-TextValue $ws.Range("E46") "  +11.38%  "
Set-TextValue $ws.Range("E47") "  +5.12%  "
Set-TextValue $ws.Range("D49") "7.82"
Set-TextValue $ws.Range("E49") "  +2.15%  "
Set-TextValue $ws.Range("D50") "0.0982"
Set-TextValue $ws.Range("E50") "  +3.78%  "
Set-TextValue $ws.Range("D51") "0.407"
Set-TextValue $ws.Range("E51") "  -0.47%  "

# Row 13/14: Polkadot and WrappedEther swap positions with updated price/volume data
Set-TextValue $ws.Range("B13") "Polkadot"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "4.22"
Set-TextValue $ws.Range("E13") "  +2.91%  "
Set-TextValue $ws.Range("B14") "WrappedEther"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D14") "1.629.66"
Set-TextValue $ws.Range("E14") "  -0.65%  "
